# Fixed fuel tax rebatements for Germany and France, storage data and German policy
$wb = $excel.ActiveWorkbook

# --- Update selections on the existing country sheets (cosmetic view state) ---
$wsDE = $wb.Worksheets.Item("DE")
$wsDE.Range("G19").Select()

$wsDK = $wb.Worksheets.Item("DK")
$wsDK.Range("G7").Select()

$wsFR = $wb.Worksheets.Item("FR")
$wsFR.Range("E9").Select()

# --- Add the new "FR_gentax" sheet right after "FR" ---
$wsGentax = $wb.Worksheets.Add($null, $wsFR)
$wsGentax.Name = "FR_gentax"

# Factor
$wsGentax.Range("B2").Value = "Factor"
$wsGentax.Range("C2").Value = 1.3

# Row labels first so the shared-string table gets HCV7/HCV8 right after "Factor"
$wsGentax.Range("B5").Value = "HCV7"
$wsGentax.Range("B6").Value = "HCV8"

# Header row
$wsGentax.Range("C4").Value = "Tax factor"
$wsGentax.Range("D4").Value = "eta"
$wsGentax.Range("E4").Value = "beta_b"
$wsGentax.Range("F4").Value = "reduction rate"
$wsGentax.Range("G4").Value = "normal tax"
$wsGentax.Range("H4").Value = "tax reduction"

# HCV7 row
$wsGentax.Range("C5").Value = 1.3
$wsGentax.Range("D5").Value = 0.89
$wsGentax.Range("E5").Value = 0.41
$wsGentax.Range("F5").Formula = "=C5*D5*(E5/(E5+1))"
$wsGentax.Range("G5").Style = "Normal 2"
$wsGentax.Range("G5").Value = 8.37
$wsGentax.Range("H5").Formula = "=-ROUND(F5*G5,2)"

# HCV8 row
$wsGentax.Range("C6").Value = 1.3
$wsGentax.Range("D6").Value = 0.89
$wsGentax.Range("E6").Value = 0.29
$wsGentax.Range("F6").Formula = "=C6*D6*(E6/(E6+1))"
$wsGentax.Range("G6").Value = 8.37
$wsGentax.Range("H6").Formula = "=-ROUND(F6*G6,2)"

$wsGentax.Range("D8").Select()
